$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table A1:E15 (header + 14 data rows for lines/extractions) gains two new
# data rows ("line7" and "line8") inserted right after "line6" (old row 7).
# This pushes the existing "extr1".."extr8" rows (old rows 8-15) down to new
# rows 10-17. We replicate that shift first (without using Range.Insert, which
# would fabricate a brand-new style combination and bloat the style table),
# then fill in the two freshly vacated rows (8 and 9) with the new line7/line8
# data, and finally fix up the one value that genuinely changed on the row
# that ends up as row 17 (in_service: TRUE -> FALSE).

# 1) Shift the old rows 8:15 down to 10:17 by copying the whole block
#    (Range.Copy(Destination) keeps the original style indices intact,
#    unlike Insert which would synthesize brand-new style combinations).
$ws.Range("A8:E15").Copy($ws.Range("A10:E17"))
$excel.CutCopyMode = $false

# 2) Populate the two new rows (8 and 9) with the line7 / line8 records.
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# 3) Fix up the A (running index) and C/D values for the shifted rows
#    (10:17) to match the target data set; A must be renumbered +2, and a
#    handful of C/D figures differ from the pre-shift source row.
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false
